$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new report rows (68 and 69) to the existing data table.

$ws.Range("A68").Value = "04c80963-c1da-46a9-beb6-4133d0481392"
$ws.Range("B68").Value = "Login with valid credentials"
$ws.Range("C68").Value = "PASSED"
$ws.Range("D68").Value = "03_28_2024_19_01_59"
$ws.Range("E68").Value = "03_28_2024_19_02_05"
$ws.Range("F68").Value = "PT5.3599447S"

$ws.Range("A69").Value = "4b04eaf8-8f23-42b1-8c34-0e8c0196e588"
$ws.Range("B69").Value = "Create Country"
$ws.Range("C69").Value = "PASSED"
$ws.Range("D69").Value = "03_28_2024_19_02_09"
$ws.Range("E69").Value = "03_28_2024_19_02_16"
$ws.Range("F69").Value = "PT7.2188974S"

# Columns D:F carry a centered style throughout the sheet; match that formatting
# for the newly appended rows.
$ws.Range("D68:F69").HorizontalAlignment = -4108
